$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header / title text (sharedStrings) ---
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Million Cubic Feet)"

# --- Insert a new row for "November" in the "Year to Date" block (old row 53) ---
$ws.Rows("53:53").Insert()

# Match formatting of the surrounding data rows (label cell + numeric cells)
$lbl = $ws.Range("A53")
$nums = $ws.Range("B53:F53")
$lbl.Borders.LineStyle = 1
$lbl.Borders.Weight = 2
$nums.Borders.LineStyle = 1
$nums.Borders.Weight = 2

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 77456
$ws.Range("C53").Value = 875
$ws.Range("D53").Value = 21756
$ws.Range("E53").Value = 3798
$ws.Range("F53").Value = 51027

# --- Update "Year to Date" annual rows (now rows 55-57 after the insert) ---
$ws.Range("A55").Value = 2014
$ws.Range("B55").Value = 790531
$ws.Range("C55").Value = 4365
$ws.Range("D55").Value = 266226
$ws.Range("E55").Value = 42751
$ws.Range("F55").Value = 477188

$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 853730
$ws.Range("C56").Value = 7317
$ws.Range("D56").Value = 259989
$ws.Range("E56").Value = 42266
$ws.Range("F56").Value = 544157

$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 876439
$ws.Range("C57").Value = 10063
$ws.Range("D57").Value = 262056
$ws.Range("E57").Value = 43852
$ws.Range("F57").Value = 560468

# --- Update the "Rolling 12 Months Ending in October" header text -> November (now row 58) ---
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# --- Update "Rolling 12 Months" annual rows (now rows 59-60 after the insert) ---
$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 928345
$ws.Range("C59").Value = 7878
$ws.Range("D59").Value = 285779
$ws.Range("E59").Value = 46150
$ws.Range("F59").Value = 588538

$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 957808
$ws.Range("C60").Value = 10806
$ws.Range("D60").Value = 285439
$ws.Range("E60").Value = 47873
$ws.Range("F60").Value = 613690
